# Insert a new data row at row 484 (pushing existing rows 484:541 down to 485:542)
# and populate the new row with the latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(484).Insert()

$ws.Range("A484").Value = 10
$ws.Range("B484").Value = "Vega Modelo de Temuco"
$ws.Range("C484").Value = "La Araucanía"
$ws.Range("D484").Value = 45212
$ws.Range("E484").Value = 9
$ws.Range("F484").Value = 100114013
$ws.Range("G484").Value = "Zanahoria"
$ws.Range("H484").Value = "Sin especificar"
$ws.Range("I484").Value = "Primera"
$ws.Range("J484").Value = 100
$ws.Range("K484").Value = 8000
$ws.Range("L484").Value = 8000
$ws.Range("M484").Value = 8000
$ws.Range("N484").Value = "`$/saco 25 kilos"
$ws.Range("O484").Value = "Región de La Araucanía"
$ws.Range("P484").Value = 320
$ws.Range("Q484").Value = 25
$ws.Range("R484").Value = "Hortaliza"

# preserve the date-number-format used by the rest of column D
$ws.Range("D484").NumberFormat = $ws.Range("D485").NumberFormat()
